# fix(database): Boolean columns issues
$wb = $excel.ActiveWorkbook

# The "client" sheet is the first sheet in the workbook
$ws = $wb.Worksheets.Item("client")

# Cell C2 currently holds the numeric value 0 (quoted-number style) which
# represents a boolean "false" value read from the database. Replace it
# with the literal text "false" (leading apostrophe forces text, avoiding
# Excel's automatic TRUE/FALSE boolean coercion) so the column is treated
# consistently as a boolean/text column.
$ws.Range("C2").Formula = "'false"

# Update the active selection on the sheet to A3 (as saved by Excel after
# the edit).
$ws.Range("A3").Select()

$wb.Save()
